# Apply "view count" (column F) updates across the sheets of the
# Hangzhou comic-convention info workbook, as produced by a fresh
# scrape/build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2847
$ws1.Range("F11").Value = 501
$ws1.Range("F13").Value = 439
$ws1.Range("F22").Value = 687
$ws1.Range("F24").Value = 159
$ws1.Range("F27").Value = 574
$ws1.Range("F28").Value = 66
$ws1.Range("F30").Value = 1716
$ws1.Range("F31").Value = 422
$ws1.Range("F33").Value = 1615
$ws1.Range("F34").Value = 232
$ws1.Range("F38").Value = 633
$ws1.Range("F42").Value = 832
$ws1.Range("F43").Value = 1533
$ws1.Range("F47").Value = 81

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 109
$ws2.Range("F9").Value = 3

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2847
$ws4.Range("F9").Value  = 501
$ws4.Range("F11").Value = 439
$ws4.Range("F17").Value = 687
$ws4.Range("F19").Value = 159
$ws4.Range("F20").Value = 109
$ws4.Range("F21").Value = 109
$ws4.Range("F26").Value = 574
$ws4.Range("F28").Value = 1716
$ws4.Range("F35").Value = 3
$ws4.Range("F43").Value = 832
$ws4.Range("F44").Value = 1533
$ws4.Range("F48").Value = 81

$wb.Save()
